$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeSchemes")
$ws.Range("C1").Value = "INFORMATIONDOMAIN"
